$wb = $excel.ActiveWorkbook

# --- 1. Update the "Status" text for row 2/3 (shared across Overview/zh-cn/de-de since it's the
#        same shared string). We only need to set it once per sheet/cell that uses it; the other
#        cells sharing the string will follow automatically, but to be safe we touch every cell
#        that shows the status text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# --- 2. Populate "Latest Target File" (F) / "Latest Handback File" (G) columns for the zh-cn
#        sheet, mirroring the source md file (A) / handoff xlf (D), plus add matching hyperlinks.
$wsZh.Range("F2").Value = "9668f9b7-f5f3-4fae-9794-50e1d46ca084.md"
$wsZh.Range("G2").Value = "9668f9b7-f5f3-4fae-9794-50e1d46ca084.be80e35b1169c0b5820dadfb212981abcc896f22.zh-cn.xlf"
$wsZh.Range("F3").Value = "a54cf248-54d7-4aca-9ddd-f24ed40d6320.md"
$wsZh.Range("G3").Value = "a54cf248-54d7-4aca-9ddd-f24ed40d6320.4a55480f235695f99e7286b65c06c0b0c7775556.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/80cfeac6b701dc2701bd6d48794d01e6a5eef93a/e2e/9668f9b7-f5f3-4fae-9794-50e1d46ca084.md", "", "", "9668f9b7-f5f3-4fae-9794-50e1d46ca084.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3de91365c034864bca3b0f80b586a21e3d2a157/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/high/9668f9b7-f5f3-4fae-9794-50e1d46ca084.be80e35b1169c0b5820dadfb212981abcc896f22.zh-cn.xlf", "", "", "9668f9b7-f5f3-4fae-9794-50e1d46ca084.be80e35b1169c0b5820dadfb212981abcc896f22.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/80cfeac6b701dc2701bd6d48794d01e6a5eef93a/e2e/a54cf248-54d7-4aca-9ddd-f24ed40d6320.md", "", "", "a54cf248-54d7-4aca-9ddd-f24ed40d6320.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a3de91365c034864bca3b0f80b586a21e3d2a157/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/high/a54cf248-54d7-4aca-9ddd-f24ed40d6320.4a55480f235695f99e7286b65c06c0b0c7775556.zh-cn.xlf", "", "", "a54cf248-54d7-4aca-9ddd-f24ed40d6320.4a55480f235695f99e7286b65c06c0b0c7775556.zh-cn.xlf")

# --- 3. Same for de-de sheet.
$wsDe.Range("F2").Value = "9668f9b7-f5f3-4fae-9794-50e1d46ca084.md"
$wsDe.Range("G2").Value = "9668f9b7-f5f3-4fae-9794-50e1d46ca084.be80e35b1169c0b5820dadfb212981abcc896f22.de-de.xlf"
$wsDe.Range("F3").Value = "a54cf248-54d7-4aca-9ddd-f24ed40d6320.md"
$wsDe.Range("G3").Value = "a54cf248-54d7-4aca-9ddd-f24ed40d6320.4a55480f235695f99e7286b65c06c0b0c7775556.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/80cfeac6b701dc2701bd6d48794d01e6a5eef93a/e2e/9668f9b7-f5f3-4fae-9794-50e1d46ca084.md", "", "", "9668f9b7-f5f3-4fae-9794-50e1d46ca084.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/951ac34663846c818454d59929511194451a71be/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/high/9668f9b7-f5f3-4fae-9794-50e1d46ca084.be80e35b1169c0b5820dadfb212981abcc896f22.de-de.xlf", "", "", "9668f9b7-f5f3-4fae-9794-50e1d46ca084.be80e35b1169c0b5820dadfb212981abcc896f22.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/80cfeac6b701dc2701bd6d48794d01e6a5eef93a/e2e/a54cf248-54d7-4aca-9ddd-f24ed40d6320.md", "", "", "a54cf248-54d7-4aca-9ddd-f24ed40d6320.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/951ac34663846c818454d59929511194451a71be/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/high/a54cf248-54d7-4aca-9ddd-f24ed40d6320.4a55480f235695f99e7286b65c06c0b0c7775556.de-de.xlf", "", "", "a54cf248-54d7-4aca-9ddd-f24ed40d6320.4a55480f235695f99e7286b65c06c0b0c7775556.de-de.xlf")

# --- 4. Latest Handback DateTime (H) updates.
$wsZh.Range("H2").Value = "2016-03-31 08:19:34"
$wsZh.Range("H3").Value = "2016-03-31 08:19:34"
$wsDe.Range("H2").Value = "2016-03-31 08:19:50"
$wsDe.Range("H3").Value = "2016-03-31 08:19:50"
